# Fruta / hortaliza, semanal
# Insert a new weekly record at row 59 (Mandarina - Murcott, Segunda) for the
# "Agrícola del Norte S.A. de Arica" market, shifting the existing rows
# 59-78 down to 60-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 59, pushing rows 59:78 down
# to 60:79 (entire-row insert, default shift = down).
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new observation.
$ws.Cells.Item(59, 1).Value  = 1
$ws.Cells.Item(59, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value  = 44524
$ws.Cells.Item(59, 5).Value  = 15
$ws.Cells.Item(59, 6).Value  = "Fruta"
$ws.Cells.Item(59, 7).Value  = 100102
$ws.Cells.Item(59, 8).Value  = "Cítricos"
$ws.Cells.Item(59, 9).Value  = 100102004
$ws.Cells.Item(59, 10).Value = "Mandarina"
$ws.Cells.Item(59, 11).Value = "Murcott"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 300
$ws.Cells.Item(59, 14).Value = 14000
$ws.Cells.Item(59, 15).Value = 15000
$ws.Cells.Item(59, 16).Value = 14500
$ws.Cells.Item(59, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(59, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(59, 19).Value = 725
$ws.Cells.Item(59, 20).Value = 20
